$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row: Date, Company, Count
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Company"
$ws.Range("C1").Value = "Count"

# Update data row 2 with new crawler volume data
$ws.Range("A2").Value = "Hour1"
$ws.Range("B2").Value = "TDM"
$ws.Range("C2").Value = 0

# Remove the old third data row entirely (was Intel/5/Week1)
$ws.Rows.Item(3).Delete()
